# "Thay doi format muc 7 cho gon hon, in ra theo mau nhu muc 12"
# The old row 3 (merged A3:AL3, blank placeholder row styled like the
# title band) is removed entirely so the header row (old row 4) becomes
# row 3 - matching the layout of the neighbouring template ("muc 12").
# Column E is also widened a bit so its header text fits on one line.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the blank styled row 3 - everything below shifts up one row,
# so the column-header row (was row 4) becomes row 3 and the
# A3:AL3 merged cell disappears along with it.
$ws.Rows(3).Delete()

# Widen column E a bit (old width 8 -> ~10.57) so the "Chuyen" header
# prints on a single line, matching the "muc 12" template.
$ws.Columns("E:E").ColumnWidth = 9.6

# Leave the selection on the new header row, matching the saved view.
$ws.Rows("3:3").Select()
